$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "React native" -> "React native (Frontend)" for the three
#        rows that used to share that technology label (rows 9-11).
$ws.Range("B9").Value = "React native (Frontend)"
$ws.Range("B10").Value = "React native (Frontend)"
$ws.Range("B11").Value = "React native (Frontend)"

# --- 2. Add the new column H ("Référance (sur le git)") by first cloning
#        the formatting of column G (same row-by-row styles) then filling
#        in the actual reference values/links.
$ws.Range("G1:G15").Copy()
$ws.Range("H1:H15").PasteSpecial(-4122)

$ws.Range("H3").Value = "Référance (sur le git)"
$ws.Range("H5").Value = "doc/Maquette app.pdf"
$ws.Range("H6").Value = "doc/Maquette app.pdf"
$ws.Range("H7").Value = "doc/Maquette app.pdf"
$ws.Range("H9").Value = "front_mobile/chibre-manager/src/screen/Home"
$ws.Range("H10").Value = "front_mobile/chibre-manager/src/screen/Game"
$ws.Range("H11").Value = "front_mobile/chibre-manager/src/screen/CreateGame"
$ws.Range("H12").Value = "doc/schéma ERD.svg"
$ws.Range("H13").Value = "backend_api/chibre-manager/db/migrate/"

# --- 3. Widen columns E and H to their new widths.
$ws.Range("E1").ColumnWidth = 13.833333333333334
$ws.Range("H1").ColumnWidth = 58.5

# --- 4. Extend the title merge A1:G2 -> A1:H2 now that there is an extra
#        column of data.
$ws.Range("A1:H2").Merge()

# --- 5. Restore the selected cell shown when the workbook was last saved.
$ws.Range("M14").Select() | Out-Null
